$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAZOPREVIR")

# Remove the "80K+156V" G3/C-WORTHY finding row (row 54) — this shifts every
# row below it up by one, per the email from Emma referenced in the commit
# message.
$ws.Rows.Item(54).Delete()

# The sheet's AutoFilter-backed _FilterDatabase name spans one row fewer now.
$wb.Names.Item("GRAZOPREVIR!_FilterDatabase").RefersTo = "=GRAZOPREVIR!`$B`$1:`$R`$185"

# Leave the cursor where the deletion happened (matches the author's saved
# selection state).
$ws.Rows.Item(31).Select()
